$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the survey response counts (Python survey results increased)
$ws.Range("C9").Value = 5
$ws.Range("C10").Value = 6
$ws.Range("C13").Value = 6

# Update the active cell selection to match the saved view state
$ws.Range("H15").Select()
